# Commit: "fix typo in sujet"
#
# The instructions paragraph reads "Nous vous demandons donc nous aider
# à refaire cette page..." but is missing the word "de": it should read
# "Nous vous demandons donc de nous aider à refaire cette page...".
#
# Locate "nous aider à" (unique in the document), collapse the found
# range to its start, and type "de " right there - exactly what typing
# the missing word into that spot in Word would do.

$d = $word.ActiveDocument

$rng = $d.Content
$found = $rng.Find.Execute(
    "nous aider à",   # old
    $true,            # MatchCase
    $false,           # MatchWholeWord
    $false,           # MatchWildcards
    $false,           # MatchSoundsLike
    $false,           # MatchAllWordForms
    $true,            # Forward
    1,                # Wrap (wdFindContinue)
    $false,           # Format
    "",                # new (unused - plain Find)
    0                 # Replace (wdReplaceNone)
)

Write-Host "Found 'nous aider à': $found"

if ($found) {
    # Collapse the found range to an insertion point right before "nous".
    $rng.End = $rng.Start
    $rng.InsertBefore("de ")
    Write-Host "Inserted missing 'de '."
}
